$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the phone numbers (and their derived "clean"/"last4" columns) for
# the four members whose numbers were re-scrubbed.
$ws.Range("D3").Value = 9356949919
$ws.Range("G3").Value = 9356949919
$ws.Range("H3").Value = 9919

$ws.Range("D4").Value = 7977603065
$ws.Range("G4").Value = 7977603065
$ws.Range("H4").Value = 3065

$ws.Range("D5").Value = 8855987935
$ws.Range("G5").Value = 8855987935
$ws.Range("H5").Value = 7935

$ws.Range("D6").Value = 9867100946
$ws.Range("G6").Value = 9867100946
$ws.Range("H6").Value = 946

# Leave the selection where the author last left it before saving.
[void]$ws.Range("H5").Select()
